# Ensure consistency of naming convention: rename ELC* indicator codes to PWR*
# on the FLO_EMIS sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FLO_EMIS")

# Rename values in the Other_Indexes / Pset_PD / Pset_CI / Cset_CN columns
$ws.Range("G4").Value = "PWRGAS"
$ws.Range("J4").Value = "PWRGAS"
$ws.Range("K4").Value = "PWRCO2N"

$ws.Range("G5").Value = "PWRGAS"
$ws.Range("J5").Value = "PWRGAS"
$ws.Range("K5").Value = "PWRCO2S"

$ws.Range("G6").Value = "PWRCOA"
$ws.Range("J6").Value = "PWRCOA"
$ws.Range("K6").Value = "PWRCO2N"

$ws.Range("G7").Value = "PWRPEA"
$ws.Range("J7").Value = "PWRPEA"
$ws.Range("K7").Value = "PWRCO2N"

$ws.Range("G8").Value = "PWRCOA"
$ws.Range("J8").Value = "PWRCOA"
$ws.Range("K8").Value = "PWRCO2S"

$ws.Range("G9").Value = "PWRPEA"
$ws.Range("J9").Value = "PWRPEA"
$ws.Range("K9").Value = "PWRCO2S"

# K9 loses its bottom border in the edit (style matches K5/K8 rather than I9)
$ws.Range("K9").Borders.Item(9).LineStyle = -4142

# Update the selection on the sheet
$ws.Activate()
$ws.Range("A3").Select()
